$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Requirements Phase Defects" --------------------------------
# The long, free-form reviewer comments in rows 11-19 (cols C:E) are removed;
# the cells stay in place (with their original styling) but are emptied out.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C11:E19").ClearContents()

# With the wrapped multi-line text gone, the explicit tall row heights are no
# longer needed - auto-fit collapses rows 11-19 back to the default height.
$ws1.Rows("11:19").EntireRow.AutoFit()

# --- Sheet 2: "Unit Design Phase Defects" ----------------------------------
# (no content changes on this sheet)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 3: "Coding Phase Defects" ---------------------------------------
# (no content changes on this sheet)
$ws3 = $wb.Worksheets.Item(3)

# --- Active tab / selection -------------------------------------------------
# Previously sheet 2 ("Unit Design Phase Defects") was the active tab with
# the selection resting on C17:E19 of... now sheet 1 becomes the active /
# selected tab, with the C17:E19 block selected (the block whose contents
# were just cleared).
$ws1.Activate()
$ws1.Range("C17:E19").Select()
